# Update docs for v0.8
# - Resources sheet gains a "Qty" column (C) with per-resource quantities
# - Resource names get padded with trailing spaces (fixed-width style)
# - Selection/active-cell bookkeeping on sheet2 and sheet3 updated

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resources")

# Add new "Qty" column next to "Cost"
$ws.Range("C1").Value = "Qty"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 1

# Pad resource names with trailing spaces to match new fixed-width look
$ws.Range("A2").Value = "Wood                                                   "
$ws.Range("A3").Value = "Metal                                                 "
$ws.Range("A4").Value = "Stone                           "

# Update the selection on the Expansion sheet
$ws2 = $wb.Worksheets.Item("Expansion")
[void]$ws2.Range("B3").Select()

# Update the selection shown on the Resources sheet, and leave it as the active tab
[void]$ws.Activate()
[void]$ws.Range("N9").Select()
